# Add "OOTB Domain Groups" worksheet after "Custom Domain Group" and
# populate it with the out-of-the-box ISP / domain-group reference table.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ootb = $wb.Sheets.Add($null, $lastSheet)
$ootb.Name = "OOTB Domain Groups"

$rows = @(
    ,@("Gmail", "gmail.com", "googlemail.com", "googlemail.co.uk")
    ,@("Microsoft", "live.com", "msn.com", "hotmail.ca", "hotmail.com", "hotmail.de", "hotmail.dk", "hotmail.co.jp", "hotmail.it", "hotmail.es", "hotmail.fr", "hotmail.co.uk", "hotmail.co.kr", "hotmail.com.au", "hotmail.com.ar", "hotmail.co.il", "hotmail.com.br", "hotmail.com.tr", "hotmail.co.th", "hotmail.jp", "hotmail.se", "live.at", "live.be", "live.ca", "live.cl", "live.cn", "live.co.kr", "live.co.uk", "live.com.ar", "live.com.au", "live.com.mx", "live.com.my", "live.com.sg", "live.de", "live.dk", "live.fr", "live.hk", "live.ie", "live.in", "live.it", "live.jp", "live.nl", "live.no", "live.ru", "live.se", "outlook.com", "live.co.uk", "hotmail.gr", "windowslive.com", "xbox.com", "hotmail.cl", "live.at", "live.jp", "live.ca", "hotmail.ca", "hotmail.ch", "live.fr", "live.it", "live.nl", "outlook.ie", "outlook.com.br", "live.com.pt", "live.be", "live.co.za", "mts.net")
    ,@("Yahoo", "yahoo.com", "rocketmail.com", "rogers.com", "sky.com", "talk21.com", "y7mail.com", "yahoo.at", "yahoo.be", "yahoo.bg", "yahoo.ca", "yahoo.cl", "yahoo.co.hu", "yahoo.co.id", "yahoo.co.il", "yahoo.co.in", "yahoo.co.jp", "yahoo.co.kr", "yahoo.com.ar", "yahoo.com.au", "yahoo.com.biz", "yahoo.com.br", "yahoo.com.cn", "yahoo.com.co", "yahoo.com.hk", "yahoo.com.hr", "yahoo.com.in", "yahoo.com.jp", "yahoo.com.kr", "yahoo.com.mx", "yahoo.com.my", "yahoo.com.net", "yahoo.com.pe", "yahoo.com.ph", "yahoo.com.sg", "yahoo.com.tr", "yahoo.com.tw", "yahoo.com.ua", "yahoo.com.ve", "yahoo.com.vn", "yahoo.co.nz", "yahoo.co.th", "yahoo.co.uk", "yahoo.co.za", "yahoo.cz", "yahoo.de", "yahoo.dk", "yahoo.ee", "yahoo.es", "yahoo.fi", "yahoo.fr", "yahoogroups.co.kr", "yahoogroups.com.cn", "yahoogroups.com.sg", "yahoogroups.com.tw", "yahoogrupper.dk", "yahoogruppi.it", "yahoo.gr", "yahoo.hr", "yahoo.hu", "yahoo.ie", "yahoo.in", "yahoo.it", "yahoo.lt", "yahoo.lv", "yahoo.nl", "yahoo.no", "yahoo.pl", "yahoo.pt", "yahoo.ro", "yahoo.rs", "yahoo.se", "yahoo.si", "yahoo.sk", "yahooxtra.co.nz", "ymail.com", "aol.com", "aim.com", "compuserve.com", "cs.com", "netscape.com", "netscape.net", "wmconnect.com", "aol.co.uk", "aol.in", "aol.de", "aol.fr", "aol.nl", "aol.pl", "aol.jp", "aol.es", "aol.it", "aol.com.ar", "aol.fi", "aol.cl", "aol.com.co", "aol.com.ve", "aol.com.au", "aol.at", "aol.be", "aol.com.br", "aol.cz", "aol.dk", "myaol.jp", "aolnorge.no", "aolpolska.pl", "aolpolcka.pl", "aolpoland.pl", "aol.ru", "aol.kr", "aol.se", "aol.ch", "aol.com.tr", "aol.co.nz", "aolchina.com", "aol.hk", "aol.tw", "luckymail.com", "verizon.net", "aol.com.mx", "bellatlantic.net", "citlink.net", "frontier.com", "frontiernet.net", "games.com", "goowy.com", "gte.net", "love.com", "verizon.net.in", "wild4music.com", "wow.com", "yahoo.cn", "yahoo.ne.jp", "yahoogroups.ca", "yahoogroups.co.in", "yahoogroups.co.uk", "yahoogroups.com", "yahoogroups.com.au", "yahoogroups.com.hk", "yahoogroups.de", "ybb.ne.jp", "ygm.com")
    ,@("Apple", "mac.com", "icloud.com", "apple.com", "me.com")
    ,@("Comcast", "comcast.net")
    ,@("Orange", "orange.fr", "orange.com", "wanadoo.fr", "francetelecom.com", "voila.fr", "voila.com")
    ,@("La Poste", "laposte.net")
    ,@("Italia Online", "libero.it", "inwind.it", "iol.it", "blu.it", "giallo.it", "virgilio.it")
    ,@("WP", "wp.pl", "o2.pl")
    ,@("United Internet", "web.de", "gmx.de", "gmx.ch", "gmx.net", "gmx.com", "gmx.at", "gmx.fr", "mail.com", "1and1.com", "1und1.de")
    ,@("Bigpond", "bigpond.com", "bigpond.net.au", "bigpond.com.au", "telstra.com", "bigpond.net")
    ,@("Docomo", "docomo.ne.jp")
    ,@("Softbank", "softbank.ne.jp", "c.vodafone.ne.jp", "d.vodafone.ne.jp", "h.vodafone.ne.jp", "k.vodafone.ne.jp", "n.vodafone.ne.jp", "q.vodafone.ne.jp", "r.vodafone.ne.jp", "s.vodafone.ne.jp", "t.vodafone.ne.jp", "jp-c.ne.jp", "jp-d.ne.jp", "jp-h.ne.jp", "jp-k.ne.jp", "jp-n.ne.jp", "jp-q.ne.jp", "jp-r.ne.jp", "jp-s.ne.jp", "jp-t.ne.jp")
    ,@("KDDI", "au.com", "ezweb.ne.jp", "uqmobile.jp")
)

$rowIndex = 1
foreach ($row in $rows) {
    $colIndex = 1
    foreach ($val in $row) {
        $ootb.Cells.Item($rowIndex, $colIndex).Value2 = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}

$ootb.Range("A23").Select()

$warmupPlan = $wb.Worksheets.Item(1)
$warmupPlan.Activate()
